$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.55%"
$ws.Range("E2").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.69%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.979"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.49%"
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.67%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.201"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.27%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.988"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.07%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9160"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.71%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09069"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-8.93%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.28%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.33%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03546"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "7.04%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09934"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.11%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001479"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.45%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005694"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.03%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.475"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.42%"
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.56%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1316"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.02%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.542"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.63%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2233"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.05%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04662"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.79%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001230"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.12%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004440"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.49%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.09%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004751"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "39.96%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01740"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.93%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04692"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.03%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007942"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.70%"
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.69%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007682"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "12.38%"
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.66%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009754"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.28%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006048"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.10%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.09%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.598"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "186.04%"
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "34.77%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.09%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.09%"
$ws.Range("E51").Style = "Normal"
